# Automatische test-sync: 2025-06-23 18:37:50
# Append a new log row (row 17) to the "Logs" sheet, extend the
# conditional formatting ranges to include it, and refresh the
# "Dashboard" category-count summary table accordingly.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 17

$logs.Cells.Item($newRow, 1).Value = "Productinformatie"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Wat is het verschil tussen product A en product B?"
$logs.Cells.Item($newRow, 4).Value = "Productinformatie"
$logs.Cells.Item($newRow, 5).Value = "Beste klant,`nDank u voor uw interesse in onze producten A en B. Product A is ontworpen met functies X en Y, terwijl product B functies P en Q biedt. Om u specifieker advies te kunnen geven, zouden we graag willen weten waarvoor u het product wilt gebruiken en welke functies voor u het belangrijkst zijn. Zo kunnen we u beter informeren over welk product het beste bij uw behoeften past.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Cells.Item($newRow, 6).Value = "2025-06-23 18:37:23"
$logs.Cells.Item($newRow, 7).Value = "Ja"

# Re-fit the row height to content (the multi-line "Antwoord" text would
# otherwise leave a stale auto-sized custom height behind).
$logs.Rows.Item($newRow).AutoFit()

# Extend the conditional-formatting ranges on columns D and G so the new
# row participates in the same "Categorie" / "Beantwoord" highlighting.
$fcCategorie = $logs.Range("D2:D16").FormatConditions
$fcCategorie.Item(1).ModifyAppliesToRange($logs.Range("D2:D17"))

$fcBeantwoord = $logs.Range("G2:G16").FormatConditions
$fcBeantwoord.Item(1).ModifyAppliesToRange($logs.Range("G2:G17"))

# Update the "Dashboard" category summary: Productinformatie now has 2
# entries and moves above "Sollicitatie / Vacature" in the ranking.
$dashboard.Cells.Item(7, 1).Value = "Productinformatie"
$dashboard.Cells.Item(7, 2).Value = 2
$dashboard.Cells.Item(8, 1).Value = "Sollicitatie / Vacature"
$dashboard.Cells.Item(8, 2).Value = 1
